$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain (non-numeric-looking) string:
# setting .Value directly keeps them as text with no style change.
$ws.Range('D2').Value = '27.216.09'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').Value = '1.772.13'
$ws.Range('E3').Value = '  +3.49%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  +1.60%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('E7').Value = '  +12.01%  '
$ws.Range('E8').Value = '  +6.66%  '
$ws.Range('E9').Value = '  +1.86%  '
$ws.Range('E10').Value = '  +0.91%  '
$ws.Range('E11').Value = '  +3.93%  '
$ws.Range('E12').Value = '  -0.02%  '
$ws.Range('E13').Value = '  +2.54%  '
$ws.Range('E14').Value = '  +3.36%  '
$ws.Range('D15').Value = '1.762.55'
$ws.Range('E15').Value = '  +3.34%  '
$ws.Range('E16').Value = '  +1.10%  '
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('E18').Value = '  +0.72%  '
$ws.Range('E19').Value = '  +0.95%  '
$ws.Range('E20').Value = '  +0.10%  '
$ws.Range('E21').Value = '  +1.49%  '
$ws.Range('E22').Value = '  +3.88%  '
$ws.Range('D23').Value = '27.290.45'
$ws.Range('E23').Value = '  +0.34%  '
$ws.Range('E24').Value = '  +4.10%  '
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('E28').Value = '  +11.30%  '
$ws.Range('D29').Value = '1.967.95'
$ws.Range('E29').Value = '  +3.50%  '
$ws.Range('E30').Value = '  +0.89%  '
$ws.Range('E31').Value = '  +4.56%  '
$ws.Range('E32').Value = '  +5.40%  '
$ws.Range('E33').Value = '  +4.77%  '
$ws.Range('E34').Value = '  +0.77%  '
$ws.Range('E35').Value = '  +1.71%  '
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('E37').Value = '  +1.38%  '
$ws.Range('B38').Value = 'TheSandbox'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E38').Value = '  +3.83%  '
$ws.Range('B39').Value = 'InternetComputer(DFINITY)'
$ws.Range('C39').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('B40').Value = 'Algorand'
$ws.Range('C40').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E40').Value = '  +0.82%  '
$ws.Range('E41').Value = '  +1.26%  '
$ws.Range('E42').Value = '  +7.52%  '
$ws.Range('E43').Value = '  +2.17%  '
$ws.Range('E44').Value = '  +3.53%  '
$ws.Range('E45').Value = '  +1.60%  '
$ws.Range('E46').Value = '  +2.12%  '
$ws.Range('E47').Value = '  +2.42%  '
$ws.Range('E48').Value = '  +2.02%  '
$ws.Range('B49').Value = 'EOS'
$ws.Range('C49').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('E49').Value = '  +2.67%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('E51').Value = '  +0.14%  '

# Cells whose new value LOOKS like a number (e.g. '1.001'): Excel's COM
# layer auto-converts such strings to floating point numbers, which loses
# the exact text / precision the source data needs. Force text storage by
# temporarily switching the cell to a text number format, assigning the
# value, then resetting the style back to 'Normal' so the cell ends up
# with the original (default) style index but a text value.
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5309'
$ws.Range('D7').Style = 'Normal'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3661'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.82'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07350'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.085'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.918'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '88.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001043'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06413'
$ws.Range('D19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.70'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.800'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.115'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.09'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.326'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.056'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09755'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.545'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.615'
$ws.Range('D34').Style = 'Normal'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02228'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05955'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.18'
$ws.Range('D37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.6130'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '4.820'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.2014'
$ws.Range('D40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.430'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.041'
$ws.Range('D42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.136'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.618'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5735'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '120.82'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.875'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.113'
$ws.Range('D49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06709'
$ws.Range('D50').Style = 'Normal'
